# Update the log entry for row 44: the end time ("C44") was missing, and the
# start time in "B44" only stored a time-of-day fraction. Fix the log so that
# B44 carries the full date/time (2017-12-02 22:00) and C44 gets the
# corresponding end time (2017-12-03 00:00, i.e. serial 43072). The dependent
# "Work Time" formula in D44 (=ABS(C44-B44)) and the Table1 totals-row formula
# in D50 (=SUM(Table1[Work Time])*24) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = 43071.916666666664
$ws.Range("C44").Value = 43072
